$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("model", "n_layers", "layer_1", "layer_2", "layer_3", "layer_4", "layer_5", "layer_6", "layer_7", "time_steps", "loss_function", "optimizer", "epochs", "validation_split", "RMSE_best", "RMSE_av")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Range("F8").Select()
